$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C. This shifts the former
# columns C, D, E one position to the right (-> D, E, F) and
# copies formatting (style) from column C into the new column.
$ws.Columns("C").Insert()

# Populate the new column C, and refresh values that were
# re-measured/changed alongside the new strategy column.
$ws.Range("B1").Value = "11100=(11111|11111)"
$ws.Range("C1").Value = "11100=(11111|11111)"
$ws.Range("D1").Value = "00111=(11011|11101)"
$ws.Range("E1").Value = "11100=(11111|11111)"
$ws.Range("F1").Value = "00111=(11011|11101)"
$ws.Range("A2").Value = "φ PyPhi"
$ws.Range("B2").Value = 0.25
$ws.Range("C2").Value = 0.25
$ws.Range("E2").Value = 0.25
$ws.Range("A3").Value = "φ Fuerza Brutal"
$ws.Range("B3").Value = 0.25
$ws.Range("C3").Value = 0.25
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0.25
$ws.Range("F3").Value = 0
$ws.Range("A4").Value = "φ Ramificación"
$ws.Range("B4").Value = 0.25
$ws.Range("C4").Value = 0.25
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.25
$ws.Range("F4").Value = 0
$ws.Range("A5").Value = "φ Genético"
$ws.Range("B5").Value = 0.25
$ws.Range("C5").Value = 0.25
$ws.Range("D5").Value = 0.25
$ws.Range("E5").Value = 0.25
$ws.Range("F5").Value = 0.25
$ws.Range("A6").Value = "═━━━━━═"
$ws.Range("B6").Value = "═━━━━━═"
$ws.Range("C6").Value = "═━━━━━═"
$ws.Range("D6").Value = "═━━━━━═"
$ws.Range("E6").Value = "═━━━━━═"
$ws.Range("F6").Value = "═━━━━━═"
$ws.Range("A7").Value = "(ms) PyPhi"
$ws.Range("B7").Value = 0.02100014686584473
$ws.Range("C7").Value = 0.02008271217346191
$ws.Range("E7").Value = 0.02014970779418945
$ws.Range("A8").Value = "(ms) Fuerza Brutal"
$ws.Range("B8").Value = 0.1882762908935547
$ws.Range("C8").Value = 0.1904587745666504
$ws.Range("D8").Value = 0.03456354141235352
$ws.Range("E8").Value = 0.1996073722839355
$ws.Range("F8").Value = 0.0331881046295166
$ws.Range("A9").Value = "(ms) Ramificación"
$ws.Range("B9").Value = 0.1299138069152832
$ws.Range("C9").Value = 0.09871888160705566
$ws.Range("D9").Value = 0.03614091873168945
$ws.Range("E9").Value = 0.09358668327331543
$ws.Range("F9").Value = 0.03635573387145996
$ws.Range("A10").Value = "(ms) Genético"
$ws.Range("B10").Value = 4.855859518051147
$ws.Range("C10").Value = 0.1205987930297852
$ws.Range("D10").Value = 0.04657077789306641
$ws.Range("E10").Value = 0.1291520595550537
$ws.Range("F10").Value = 0.04761123657226562
$ws.Range("A11").Value = "═━━━━━═"
$ws.Range("B11").Value = "═━━━━━═"
$ws.Range("C11").Value = "═━━━━━═"
$ws.Range("D11").Value = "═━━━━━═"
$ws.Range("E11").Value = "═━━━━━═"
$ws.Range("F11").Value = "═━━━━━═"
